$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing columns C, D, E (old Status / IPv4 / Created At data)
$ws.Range("C1:E2").Clear()

# Set header row
$ws.Range("A1").Value = "Container Name"
$ws.Range("B1").Value = "Image Name"

# Set data row
$ws.Range("A2").Value = "test-import-container"
$ws.Range("B2").Value = "nginx:stable-alpine-perl"

# Column widths (closest achievable values given engine's pixel-grid snapping;
# targets are 21.290714285714284 / 23.290714285714284 chars-equivalent)
$ws.Columns.Item(1).ColumnWidth = 20.5
$ws.Columns.Item(2).ColumnWidth = 22.5

# Row heights
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75
